$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 22499.4
$ws.Range("I40").Value = 16785.715
$ws.Range("J40").Value = 35831.332
$ws.Range("K40").Value = 16785.715
$ws.Range("L40").Value = 35831.332
$ws.Range("M40").Value = -16610.715
$ws.Range("N40").Value = -36181.332
$ws.Range("H54").Value = 15000
$ws.Range("I54").Value = 11666.667
$ws.Range("K54").Value = 11666.667
$ws.Range("M54").Value = -11180.667
$ws.Range("H69").Value = 83349464
$ws.Range("J69").Value = 21500
$ws.Range("L69").Value = 64500
$ws.Range("N69").Value = -66248
$ws.Range("H72").Value = 83349464
$ws.Range("J72").Value = 21500
$ws.Range("L72").Value = 193500
$ws.Range("N72").Value = -202236
$ws.Range("H74").Value = 3861.7896
$ws.Range("I74").Value = 3834
$ws.Range("J74").Value = 3900
$ws.Range("K74").Value = 3834
$ws.Range("L74").Value = 3900
$ws.Range("M74").Value = -2898
$ws.Range("N74").Value = -5772
$ws.Range("H76").Value = 5361.875
$ws.Range("I76").Value = 3950
$ws.Range("J76").Value = 5832.5
$ws.Range("K76").Value = 3950
$ws.Range("L76").Value = 5832.5
$ws.Range("M76").Value = -3635
$ws.Range("N76").Value = -6462.5
$ws.Range("H77").Value = 3861.7896
$ws.Range("I77").Value = 3834
$ws.Range("J77").Value = 3900
$ws.Range("K77").Value = 19170
$ws.Range("L77").Value = 19500
$ws.Range("M77").Value = -14490
$ws.Range("N77").Value = -28860
$ws.Range("H79").Value = 5361.875
$ws.Range("I79").Value = 3950
$ws.Range("J79").Value = 5832.5
$ws.Range("K79").Value = 3950
$ws.Range("L79").Value = 5832.5
$ws.Range("M79").Value = -2858
$ws.Range("N79").Value = -8016.5
$ws.Range("H92").Value = 83333800
$ws.Range("I92").Value = 111111620
$ws.Range("K92").Value = 111111620
$ws.Range("M92").Value = -111110372
$ws.Range("H100").Value = 711
$ws.Range("I100").Value = 718.25
$ws.Range("K100").Value = 718.25
$ws.Range("M100").Value = -177.25
$ws.Range("H132").Value = 4268.769
$ws.Range("I132").Value = 1193.75
$ws.Range("J132").Value = 12096.091
$ws.Range("K132").Value = 3581.25
$ws.Range("L132").Value = 36288.273
$ws.Range("M132").Value = -1051.25
$ws.Range("N132").Value = -41348.273
$ws.Range("H137").Value = 11448244
$ws.Range("I137").Value = 670459.6
$ws.Range("J137").Value = 22226028
$ws.Range("K137").Value = 2011378.8
$ws.Range("L137").Value = 66678084
$ws.Range("M137").Value = -2008828.8
$ws.Range("N137").Value = -66683184

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 960.65216
$ws.Range("I2").Value = 895.2273
$ws.Range("K2").Value = 895.2273
$ws.Range("M2").Value = -782.2273
$ws.Range("H32").Value = 11597.186
$ws.Range("I32").Value = 10926.896
$ws.Range("K32").Value = 10926.896
$ws.Range("M32").Value = -10639.896
$ws.Range("H61").Value = 3833.8
$ws.Range("I61").Value = 2383.2173
$ws.Range("K61").Value = 2383.2173
$ws.Range("M61").Value = -2171.2173
$ws.Range("H96").Value = 25172
$ws.Range("J96").Value = 25172
$ws.Range("L96").Value = 25172
$ws.Range("N96").Value = -30664
$ws.Range("H102").Value = 816.8484999999999
$ws.Range("I102").Value = 719.6539
$ws.Range("J102").Value = 1177.8572
$ws.Range("K102").Value = 719.6539
$ws.Range("L102").Value = 1177.8572
$ws.Range("M102").Value = 902.3461
$ws.Range("N102").Value = -4421.8572
$ws.Range("H116").Value = 960.65216
$ws.Range("I116").Value = 895.2273
$ws.Range("K116").Value = 895.2273
$ws.Range("M116").Value = 1398.7727
$ws.Range("H132").Value = 18111.578
$ws.Range("I132").Value = 21154.666
$ws.Range("K132").Value = 63463.99800000001
$ws.Range("M132").Value = -60933.99800000001
$ws.Range("H136").Value = 3833.8
$ws.Range("I136").Value = 2383.2173
$ws.Range("K136").Value = 7149.651899999999
$ws.Range("M136").Value = -4599.651899999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 960.65216
$ws.Range("I3").Value = 895.2273
$ws.Range("K3").Value = 895.2273
$ws.Range("M3").Value = -781.2273
$ws.Range("H80").Value = 867.4211
$ws.Range("I80").Value = 428.5
$ws.Range("J80").Value = 1070
$ws.Range("K80").Value = 428.5
$ws.Range("L80").Value = 1070
$ws.Range("M80").Value = 569.5
$ws.Range("N80").Value = -3066
$ws.Range("H83").Value = 867.4211
$ws.Range("I83").Value = 428.5
$ws.Range("J83").Value = 1070
$ws.Range("K83").Value = 2142.5
$ws.Range("L83").Value = 5350
$ws.Range("M83").Value = 2849.5
$ws.Range("N83").Value = -15334
$ws.Range("H86").Value = 3709.5454
$ws.Range("I86").Value = 3180.5
$ws.Range("K86").Value = 3180.5
$ws.Range("M86").Value = -2057.5
$ws.Range("H89").Value = 3709.5454
$ws.Range("I89").Value = 3180.5
$ws.Range("K89").Value = 15902.5
$ws.Range("M89").Value = -10286.5
$ws.Range("H107").Value = 2366.6758
$ws.Range("I107").Value = 2148.08
$ws.Range("K107").Value = 2148.08
$ws.Range("M107").Value = -228.0799999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 1200
$ws.Range("J14").Value = 1200
$ws.Range("L14").Value = 1200
$ws.Range("N14").Value = -1540
$ws.Range("H107").Value = 807
$ws.Range("I107").Value = 533.1111
$ws.Range("K107").Value = 533.1111
$ws.Range("M107").Value = 1386.8889
$ws.Range("H141").Value = 166875.05
$ws.Range("J141").Value = 204048.23
$ws.Range("L141").Value = 204048.23
$ws.Range("N141").Value = -214408.23

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 91037270
$ws.Range("J37").Value = 91037270
$ws.Range("L37").Value = 273111810
$ws.Range("N37").Value = -273112034
$ws.Range("H131").Value = 12347972
$ws.Range("J131").Value = 13335671
$ws.Range("L131").Value = 40007013
$ws.Range("N131").Value = -40017093

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 290.30234
$ws.Range("I2").Value = 209.09525
$ws.Range("K2").Value = 209.09525
$ws.Range("M2").Value = -96.09524999999999
$ws.Range("H26").Value = 20000
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H50").Value = 20000
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H80").Value = 58373.24
$ws.Range("I80").Value = 62519.94
$ws.Range("K80").Value = 62519.94
$ws.Range("M80").Value = -61521.94
$ws.Range("H83").Value = 58373.24
$ws.Range("I83").Value = 62519.94
$ws.Range("K83").Value = 312599.7
$ws.Range("M83").Value = -307607.7
$ws.Range("H92").Value = 15849.8
$ws.Range("J92").Value = 15849.8
$ws.Range("L92").Value = 15849.8
$ws.Range("N92").Value = -19593.8
$ws.Range("H93").Value = 44000
$ws.Range("J93").Value = 44000
$ws.Range("L93").Value = 44000
$ws.Range("N93").Value = -47744
$ws.Range("H97").Value = 1801.909
$ws.Range("I97").Value = 1091.6666
$ws.Range("K97").Value = 1091.6666
$ws.Range("M97").Value = -595.6666
$ws.Range("H126").Value = 4301.3887
$ws.Range("I126").Value = 4835.7856
$ws.Range("J126").Value = 2431
$ws.Range("K126").Value = 14507.3568
$ws.Range("L126").Value = 7293
$ws.Range("M126").Value = -12037.3568
$ws.Range("N126").Value = -12233

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 2399.3333
$ws.Range("J31").Value = 3375
$ws.Range("L31").Value = 3375
$ws.Range("N31").Value = -3871
$ws.Range("H40").Value = 17640228
$ws.Range("I40").Value = 8931463
$ws.Range("K40").Value = 8931463
$ws.Range("M40").Value = -8931327
$ws.Range("H93").Value = 2062.5293
$ws.Range("I93").Value = 1823.909
$ws.Range("K93").Value = 1823.909
$ws.Range("M93").Value = -575.9090000000001
$ws.Range("H122").Value = 10872
$ws.Range("I122").Value = 4950
$ws.Range("J122").Value = 11661.6
$ws.Range("K122").Value = 14850
$ws.Range("L122").Value = 34984.8
$ws.Range("M122").Value = -12400
$ws.Range("N122").Value = -39884.8
$ws.Range("H136").Value = 4642.559
$ws.Range("I136").Value = 2844.7144
$ws.Range("K136").Value = 8534.143199999999
$ws.Range("M136").Value = -5984.143199999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1652
$ws.Range("I107").Value = 1718.8125
$ws.Range("K107").Value = 5156.4375
$ws.Range("M107").Value = -3236.4375
$ws.Range("H122").Value = 4343.7646
$ws.Range("I122").Value = 2000.3334
$ws.Range("J122").Value = 4845.9287
$ws.Range("K122").Value = 6001.0002
$ws.Range("L122").Value = 14537.7861
$ws.Range("M122").Value = -3551.0002
$ws.Range("N122").Value = -19437.7861
$ws.Range("H136").Value = 6945.8887
$ws.Range("I136").Value = 5835.857
$ws.Range("J136").Value = 7652.273
$ws.Range("K136").Value = 17507.571
$ws.Range("L136").Value = 22956.819
$ws.Range("M136").Value = -14957.571
$ws.Range("N136").Value = -28056.819

